$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.565.29'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '1.695.34'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').Value = "'" + '1.007'
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = "'" + '309.91'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = "'" + '1.004'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').Value = "'" + '0.3728'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +2.33%  '
$ws.Range('D9').Value = "'" + '0.3421'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = "'" + '1.170'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').Value = "'" + '0.07417'
$ws.Range('E11').Value = '  +1.52%  '
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = "'" + '20.61'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = "'" + '6.200'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = "'" + '6.874'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').Value = '1.694.30'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').Value = "'" + '0.00001112'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = "'" + '1.004'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('D19').Value = "'" + '0.06704'
$ws.Range('D20').Value = "'" + '82.69'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').Value = "'" + '16.96'
$ws.Range('E21').Value = '  +2.95%  '
$ws.Range('D22').Value = "'" + '6.302'
$ws.Range('E22').Value = '  +2.51%  '
$ws.Range('D23').Value = "'" + '12.76'
$ws.Range('E23').Value = '  +6.24%  '
$ws.Range('D24').Value = '24.523.13'
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').Value = "'" + '2.446'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').Value = "'" + '2.734'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').Value = "'" + '20.12'
$ws.Range('E27').Value = '  +3.00%  '
$ws.Range('D28').Value = "'" + '149.59'
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'" + '130.24'
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('B30').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C30').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D30').Value = '1.883.32'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').Value = "'" + '1.154'
$ws.Range('E31').Value = '  +15.67%  '
$ws.Range('D32').Value = "'" + '6.555'
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('D33').Value = "'" + '4.232'
$ws.Range('E33').Value = '  +2.50%  '
$ws.Range('D34').Value = "'" + '1.770'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = "'" + '0.08726'
$ws.Range('E35').Value = '  +2.41%  '
$ws.Range('D36').Value = "'" + '13.37'
$ws.Range('E36').Value = '  +6.35%  '
$ws.Range('D37').Value = "'" + '5.485'
$ws.Range('E37').Value = '  +1.93%  '
$ws.Range('D38').Value = "'" + '0.06461'
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('D39').Value = "'" + '0.02361'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').Value = "'" + '8.826'
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').Value = "'" + '1.267'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').Value = "'" + '0.6345'
$ws.Range('E43').Value = '  +2.49%  '
$ws.Range('D44').Value = "'" + '1.006'
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('D45').Value = "'" + '13.70'
$ws.Range('E45').Value = '  +3.19%  '
$ws.Range('D46').Value = "'" + '3.810'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').Value = "'" + '0.6019'
$ws.Range('E47').Value = '  +0.96%  '
$ws.Range('D48').Value = "'" + '2.096'
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('D49').Value = "'" + '128.05'
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('D51').Value = "'" + '78.68'
$ws.Range('E51').Value = '  +2.50%  '
